$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 319-320, shifting the old 319-327 block down to 321-329
$ws.Range("A319:A320").EntireRow.Insert()

# Row 319: new weekly entry
$ws.Range("A319").Value = 11
$ws.Range("B319").Value = "Vega Monumental Concepción"
$ws.Range("C319").Value = "Bíobío"
$ws.Range("D319").Value = 44939
$ws.Range("E319").Value = 8
$ws.Range("F319").Value = 100112045
$ws.Range("G319").Value = "Zapallo"
$ws.Range("H319").Value = "Camote"
$ws.Range("I319").Value = "1a (cosecha)"
$ws.Range("J319").Value = 300
$ws.Range("K319").Value = 600
$ws.Range("L319").Value = 600
$ws.Range("M319").Value = 600
$ws.Range("N319").Value = "$/kilo (volumen en unidades)"
$ws.Range("O319").Value = "Región de O'Higgins"
$ws.Range("P319").Value = 600
$ws.Range("Q319").Value = 1
$ws.Range("R319").Value = "Hortaliza"

# Row 320: new weekly entry
$ws.Range("A320").Value = 11
$ws.Range("B320").Value = "Vega Monumental Concepción"
$ws.Range("C320").Value = "Bíobío"
$ws.Range("D320").Value = 44939
$ws.Range("E320").Value = 8
$ws.Range("F320").Value = 100112045
$ws.Range("G320").Value = "Zapallo"
$ws.Range("H320").Value = "Camote"
$ws.Range("I320").Value = "2a (cosecha)"
$ws.Range("J320").Value = 300
$ws.Range("K320").Value = 500
$ws.Range("L320").Value = 500
$ws.Range("M320").Value = 500
$ws.Range("N320").Value = "$/kilo (volumen en unidades)"
$ws.Range("O320").Value = "Región de O'Higgins"
$ws.Range("P320").Value = 500
$ws.Range("Q320").Value = 1
$ws.Range("R320").Value = "Hortaliza"
